# Update "paises.xlsx" country/provincia data and the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Update the "Datos actualizados..." timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 16:22"

# 2) Row 4 - Estados Unidos
$ws.Range("B4").Value = 888881
$ws.Range("C4").Value = 2439
$ws.Range("D4").Value = 89221
$ws.Range("E4").Value = 749291
$ws.Range("F4").Value = 15029
$ws.Range("G4").Value = 133
$ws.Range("H4").Value = 50369

# 3) Row 8 - Reino Unido
$ws.Range("B8").Value = 153584
$ws.Range("C8").Value = 455
$ws.Range("E8").Value = 41207

# 4) Row 14 - China
$ws.Range("B14").Value = 50512
$ws.Range("C14").Value = 1020
$ws.Range("E14").Value = 20574
$ws.Range("G14").Value = 52
$ws.Range("H14").Value = 3365

# 5) Row 23
$ws.Range("F23").Value = 547

# 6) Row 52
$ws.Range("D52").Value = 2500
$ws.Range("E52").Value = 1718

# 7) Row 81
$ws.Range("B81").Value = 1326
$ws.Range("C81").Value = 26
$ws.Range("D81").Value = 337
$ws.Range("E81").Value = 932
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 57

# 8) Row 104
$ws.Range("E104").Value = 276
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 7

# 9) Row 109
$ws.Range("D109").Value = 132
$ws.Range("E109").Value = 294
